$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 106: STDEV.S formulas for B:F over rows 2-104.
# B106 is entered individually (its own formula, not part of a fill group),
# while C106:F106 is entered as one range-fill so Excel groups it as a
# shared formula anchored at C106 (matching the prior AVERAGE row's pattern).
$ws.Range("B106").Formula = "=STDEV.S(B2:B104)"
$ws.Range("C106:F106").Formula = "=STDEV.S(C2:C104)"

# Row 107: margin-of-error formulas referencing row 106, same fill pattern.
$ws.Range("B107").Formula = "=B106/SQRT(103)*1.96"
$ws.Range("C107:F107").Formula = "=C106/SQRT(103)*1.96"

# Update selection to match the post-edit state
$ws.Range("H109").Select()
